$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell $ws "D2" "328.60"
Set-TextCell $ws "E2" "0.06%"
Set-TextCell $ws "G2" "16"

Set-TextCell $ws "D3" "44.05"
Set-TextCell $ws "E3" "-0.05%"
Set-TextCell $ws "G3" "16"

Set-TextCell $ws "D4" "5.495"
Set-TextCell $ws "E4" "0.16%"
Set-TextCell $ws "G4" "16"

Set-TextCell $ws "D5" "0.08358"
Set-TextCell $ws "E5" "4.10%"
Set-TextCell $ws "G5" "16"

Set-TextCell $ws "D6" "2.010"
Set-TextCell $ws "E6" "-0.89%"
Set-TextCell $ws "G6" "16"

Set-TextCell $ws "D7" "0.9753"
Set-TextCell $ws "E7" "2.29%"
Set-TextCell $ws "G7" "16"

Set-TextCell $ws "D8" "2.544"
Set-TextCell $ws "E8" "-2.63%"
Set-TextCell $ws "G8" "16"

Set-TextCell $ws "D9" "0.1102"
Set-TextCell $ws "E9" "0.28%"
Set-TextCell $ws "G9" "16"

Set-TextCell $ws "D10" "0.1907"
Set-TextCell $ws "E10" "1.15%"
Set-TextCell $ws "G10" "16"

Set-TextCell $ws "D11" "9.406"
Set-TextCell $ws "E11" "-7.04%"
Set-TextCell $ws "G11" "16"

Set-TextCell $ws "D12" "0.09682"
Set-TextCell $ws "E12" "-3.11%"
Set-TextCell $ws "G12" "16"

Set-TextCell $ws "D13" "0.04699"
Set-TextCell $ws "E13" "-0.83%"
Set-TextCell $ws "G13" "16"

Set-TextCell $ws "E14" "0.04%"
Set-TextCell $ws "G14" "16"

Set-TextCell $ws "D15" "0.001282"
Set-TextCell $ws "E15" "1.23%"
Set-TextCell $ws "G15" "16"

Set-TextCell $ws "D16" "0.005928"
Set-TextCell $ws "E16" "0.41%"
Set-TextCell $ws "G16" "16"

Set-TextCell $ws "D17" "3.389"
Set-TextCell $ws "E17" "0.46%"
Set-TextCell $ws "G17" "16"

Set-TextCell $ws "D18" "4.456"
Set-TextCell $ws "E18" "0.94%"
Set-TextCell $ws "G18" "16"

Set-TextCell $ws "D19" "0.3292"
Set-TextCell $ws "E19" "-3.55%"
Set-TextCell $ws "G19" "16"

Set-TextCell $ws "D20" "0.1374"
Set-TextCell $ws "E20" "-1.90%"
Set-TextCell $ws "G20" "16"

Set-TextCell $ws "D21" "0.2557"
Set-TextCell $ws "E21" "-0.96%"
Set-TextCell $ws "G21" "16"

Set-TextCell $ws "E22" "2.57%"
Set-TextCell $ws "G22" "16"

Set-TextCell $ws "D23" "0.001302"
Set-TextCell $ws "E23" "-0.54%"
Set-TextCell $ws "G23" "16"

Set-TextCell $ws "D24" "0.004438"
Set-TextCell $ws "E24" "2.29%"
Set-TextCell $ws "G24" "16"

Set-TextCell $ws "D25" "0.0001306"
Set-TextCell $ws "E25" "4.28%"
Set-TextCell $ws "G25" "16"

Set-TextCell $ws "G26" "16"

Set-TextCell $ws "G27" "16"

Set-TextCell $ws "G28" "16"

Set-TextCell $ws "G29" "16"

Set-TextCell $ws "G30" "16"

Set-TextCell $ws "G31" "16"

Set-TextCell $ws "G32" "16"

Set-TextCell $ws "G33" "16"

Set-TextCell $ws "G34" "16"

Set-TextCell $ws "G35" "16"

Set-TextCell $ws "G36" "16"

Set-TextCell $ws "G37" "16"

Set-TextCell $ws "D38" "0.02720"
Set-TextCell $ws "E38" "5.65%"
Set-TextCell $ws "G38" "16"

Set-TextCell $ws "D39" "0.05638"
Set-TextCell $ws "E39" "-0.57%"
Set-TextCell $ws "G39" "16"

Set-TextCell $ws "D40" "0.007853"
Set-TextCell $ws "E40" "1.27%"
Set-TextCell $ws "G40" "16"

Set-TextCell $ws "E41" "1.64%"
Set-TextCell $ws "G41" "16"

Set-TextCell $ws "D42" "0.007396"
Set-TextCell $ws "E42" "0.38%"
Set-TextCell $ws "G42" "16"

Set-TextCell $ws "D43" "0.002123"
Set-TextCell $ws "E43" "5.55%"
Set-TextCell $ws "G43" "16"

Set-TextCell $ws "D44" "0.008624"
Set-TextCell $ws "E44" "1.14%"
Set-TextCell $ws "G44" "16"

Set-TextCell $ws "D45" "0.3370"
Set-TextCell $ws "G45" "16"

Set-TextCell $ws "D46" "0.00006883"
Set-TextCell $ws "E46" "-2.99%"
Set-TextCell $ws "G46" "16"

Set-TextCell $ws "D47" "0.00000000753"
Set-TextCell $ws "E47" "0.27%"
Set-TextCell $ws "G47" "16"

Set-TextCell $ws "E48" "0.13%"
Set-TextCell $ws "G48" "16"

Set-TextCell $ws "D49" "0.003517"
Set-TextCell $ws "E49" "0.31%"
Set-TextCell $ws "G49" "16"

Set-TextCell $ws "D50" "0.003546"
Set-TextCell $ws "E50" "1.20%"
Set-TextCell $ws "G50" "16"

Set-TextCell $ws "D51" "0.00002110"
Set-TextCell $ws "E51" "0.27%"
Set-TextCell $ws "G51" "16"
